$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list price/volume refresh (GitHub Actions bot).
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the source inlineStr cells, preserving exact
# formatting such as trailing zeros and multi-dot "thousands" prices),
# then the style is reset to 'Normal' so no stray number-format /
# quote-prefix styling is left behind on the cell.

$ws.Range('D2').Value = "'27.304.99"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +1.01%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.856.53"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +1.48%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D5').Value = "'314.83"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.69%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D7').Value = "'0.4606"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.31%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.3711"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.21%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.07300"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.47%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.8899"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +1.71%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'20.09"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +1.43%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.07829"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.46%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'1.824.32"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +1.00%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'5.393"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.09%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'6.529"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -0.76%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'91.51"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -0.14%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  -0.29%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'0.000008924"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.12%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.49%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'14.76"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.36%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'27.318.07"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +0.21%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'5.121"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.35%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'10.56"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +0.06%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'2.070.77"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +3.07%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'1.920"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +3.92%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'152.19"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.64%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'18.46"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.38%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'2.055"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +0.18%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'116.07"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.74%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'5.072"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -1.24%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'0.08830"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.40%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'0.7730"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +5.64%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'3.086"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +4.42%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'1.174"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +3.46%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'4.514"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +1.50%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'2.731"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +12.23%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.37%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.01956"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +0.70%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +0.66%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'2.960"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.55%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'7.045"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -1.28%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.5131"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.39%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.1641"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.60%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'8.426"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +2.46%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.4794"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.74%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'10.37"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +2.01%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'1.003"
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Value = "'102.92"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.37%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'1.644"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +0.99%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.06227"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.22%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'65.83"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +2.20%  "
$ws.Range('E51').Style = 'Normal'
